$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zestawienie kosztów")
$pt = $ws.PivotTables(1)
$pf = $pt.PivotFields("Typ")
$df = $pt.PivotFields("Suma z Kwota")
Write-Output "DataField name: $($df.Name)"
$pf.PivotItems().Item("Koszty kredytu").Visible = $true
$pt.RefreshTable()
try {
  $pf.AutoSort(2, $df.Name)
  Write-Output "ok1"
} catch { Write-Output "err1 $_" }
$pt.RefreshTable()
